$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Hunk 1: merge the "n" / "-1" runs inside a_{n-1}'s subscript into a single
# run with text "n-1" (OMath #2 : A={a0, a1, a2, ..., a_{n-1}})
# ---------------------------------------------------------------------------
$om1 = $d.OMaths.Item(2)
$xml1 = '<m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t>A={</m:t></m:r><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t>a</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t>0</m:t></m:r></m:sub></m:sSub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t xml:space="preserve">, </m:t></m:r><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t>a</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t>1</m:t></m:r></m:sub></m:sSub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t xml:space="preserve">, </m:t></m:r><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t>a</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t>2</m:t></m:r></m:sub></m:sSub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t xml:space="preserve">, …, </m:t></m:r><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t>a</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t>n-1</m:t></m:r></m:sub></m:sSub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman" w:hint="eastAsia"/></w:rPr><m:t>}</m:t></m:r></m:oMath>'
$om1.Range.InsertXML($xml1)

# ---------------------------------------------------------------------------
# Hunk 2: drop the stray _GoBack bookmark from a_4's subscript
# (OMath #25 : [0, 0, 1, 1, 1] -> {a2, a3, a4})
# ---------------------------------------------------------------------------
$om2 = $d.OMaths.Item(25)
$xml2 = '<m:oMathPara><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman" w:hint="eastAsia"/></w:rPr><m:t>[0</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t xml:space="preserve">, </m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman" w:hint="eastAsia"/></w:rPr><m:t>0</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t xml:space="preserve">, 1, </m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman" w:hint="eastAsia"/></w:rPr><m:t>1</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t xml:space="preserve">, </m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman" w:hint="eastAsia"/></w:rPr><m:t>1]</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t xml:space="preserve"> -&gt; {</m:t></m:r><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t>a</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman" w:hint="eastAsia"/></w:rPr><m:t>2</m:t></m:r></m:sub></m:sSub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t xml:space="preserve">, </m:t></m:r><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t>a</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman" w:hint="eastAsia"/></w:rPr><m:t>3</m:t></m:r></m:sub></m:sSub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t xml:space="preserve">, </m:t></m:r><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t>a</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman" w:hint="eastAsia"/></w:rPr><m:t>4</m:t></m:r></m:sub></m:sSub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t>}</m:t></m:r></m:oMath></m:oMathPara>'
$om2.Range.InsertXML($xml2)

# ---------------------------------------------------------------------------
# Hunk 3: rewrite the "time complexity" equation from O(n^2) to
# C_m^n = n! / (m!(n-m)!), carrying the _GoBack bookmark along with it.
# (OMath #41 : O(n^2))
# ---------------------------------------------------------------------------
$om3 = $d.OMaths.Item(41)
$xml3 = '<m:oMath><m:sSubSup><m:sSubSupPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr></m:ctrlPr></m:sSubSupPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t>C</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t>m</m:t></m:r></m:sub><m:sup><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t>n</m:t></m:r></m:sup></m:sSubSup><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t>=</m:t></m:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><m:f><m:fPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/><w:i/></w:rPr></m:ctrlPr></m:fPr><m:num><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t>n!</m:t></m:r></m:num><m:den><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t>m!</m:t></m:r><m:d><m:dPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/><w:i/></w:rPr></m:ctrlPr></m:dPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t>n-m</m:t></m:r></m:e></m:d><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="宋体" w:hAnsi="Cambria Math" w:cs="Times New Roman"/></w:rPr><m:t>!</m:t></m:r></m:den></m:f></m:oMath>'
$om3.Range.InsertXML($xml3)

Write-Output "done"
